# izradena funkcija za postavlanje 'nisu naplaceni do', 'datum od' i 'datum do'
# za uvezivanje excel tablice i kreiranje novog obrasca preko aplikacije

$wb = $excel.ActiveWorkbook

# --- "Zaglavlje" sheet: update "nisu naplaceni do" date and move selection ---
$wsZaglavlje = $wb.Worksheets.Item("Zaglavlje")
$wsZaglavlje.Range("G2").Value = 42460
$wsZaglavlje.Range("G3").Select()

# --- "Racuni" sheet: update 'datum od' / 'datum do' dates, move selection ---
$wsRacuni = $wb.Worksheets.Item("Racuni")
$wsRacuni.Range("B3").Value = 42439
$wsRacuni.Range("C3").Value = 42439

# Make "Racuni" the active tab/sheet (was "Zaglavlje" before the edit)
$wsRacuni.Activate()
$wsRacuni.Range("D5").Select()
